# "Evaluation Protocol once again x2"
# Re-score a handful of the JS-SPA self-evaluation rows and leave the
# selection parked near the bottom of the sheet (where the last edit -
# the Total Score formula - lives).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numbers of Commits in GitHub: 22 -> 25
$ws.Range("C9").Value = 25

# AngularJS Project Structure: 3 -> 4
$ws.Range("C12").Value = 4

# Register Screen: 5 -> 10
$ws.Range("C17").Value = 10

# Authorization Checks: 2 -> 3
$ws.Range("C32").Value = 3

# C51 (=SUM(C6:C50)) recalculates automatically: 98 -> 108

# Leave the view scrolled down / selection on the last touched cell.
$ws.Range("E46").Select()
